$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update the (shared) product name text - both sheets reference the same
# shared string, so setting the same new text on both keeps them pointed
# at a single, updated shared-string entry.
$newProductName = "4201-RBI-EI-DB-DL-REC-INT-RNI-FFC-SAR-FFROP-DAILY-1-CTRFD-MD-TR-1-ONTIME-PER-1st"
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Change the short name from the numeric 4201 to the text "420f"
$ws1.Range("B2").Value = "420f"

# Removed test-case inter-dependency: reset each sheet's selection/scroll
# back to B1 (instead of the previously-saved mid-sheet selection).
$ws2.Activate()
$ws2.Range("B1").Select()
$ws1.Activate()
$ws1.Range("B1").Select()
